# Update the cryptocurrency price/volume list on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25 / 26 swap: Monero and Cosmos traded places in the ranking ---
# Row 25 becomes Cosmos, Row 26 becomes Monero (with refreshed data).
$ws.Cells.Item(25, 2).Value = "Cosmos"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "9.206"
$ws.Cells.Item(25, 5).Value = "  -0.63%  "

$ws.Cells.Item(26, 2).Value = "Monero"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "165.26"
$ws.Cells.Item(26, 5).Value = "  +0.96%  "

# --- Price (column D) and Volume(1h) (column E) refreshes for all other rows ---
# For price values that look like plain numbers (single '.' decimal point),
# the cell is first forced to Text format ("@") so Excel keeps storing the
# refreshed price as text instead of silently converting it to a number --
# matching how these cells were already stored as text in the source file.

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "30.617.76"
$ws.Cells.Item(2, 5).Value = "  +0.82%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "1.878.16"
$ws.Cells.Item(3, 5).Value = "  -0.06%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 5).Value = "  -0.03%  "

# Row 5 - BNB
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "248.27"
$ws.Cells.Item(5, 5).Value = "  +1.61%  "

# Row 6 - USDC
$ws.Cells.Item(6, 5).Value = "  +0.01%  "

# Row 7 - XRP
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4764"
$ws.Cells.Item(7, 5).Value = "  -0.20%  "

# Row 8 - Cardano
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2928"
$ws.Cells.Item(8, 5).Value = "  +1.54%  "

# Row 9 - Dogecoin
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06530"
$ws.Cells.Item(9, 5).Value = "  +0.11%  "

# Row 10 - Solana
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "22.06"
$ws.Cells.Item(10, 5).Value = "  +3.19%  "

# Row 11 - TRON
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07737"
$ws.Cells.Item(11, 5).Value = "  -0.32%  "

# Row 12 - Polygon
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.7415"
$ws.Cells.Item(12, 5).Value = "  +0.78%  "

# Row 13 - Litecoin
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "96.82"
$ws.Cells.Item(13, 5).Value = "  +0.29%  "

# Row 14 - WrappedEther
$ws.Cells.Item(14, 4).Value = "1.874.81"
$ws.Cells.Item(14, 5).Value = "  -0.23%  "

# Row 15 - Polkadot
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "5.204"
$ws.Cells.Item(15, 5).Value = "  +1.43%  "

# Row 16 - BitcoinCash
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "274.60"
$ws.Cells.Item(16, 5).Value = "  -0.48%  "

# Row 17 - WrappedBTC
$ws.Cells.Item(17, 4).Value = "30.716.80"
$ws.Cells.Item(17, 5).Value = "  +1.16%  "

# Row 18 - Avalanche
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "13.26"
$ws.Cells.Item(18, 5).Value = "  -0.98%  "

# Row 19 - ShibaInu
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.000007538"
$ws.Cells.Item(19, 5).Value = "  -0.18%  "

# Row 20 - Dai
$ws.Cells.Item(20, 5).Value = "  +0.02%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Cells.Item(21, 4).Value = "2.123.17"
$ws.Cells.Item(21, 5).Value = "  -0.10%  "

# Row 22 - BinanceUSD
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.9995"
$ws.Cells.Item(22, 5).Value = "  -0.07%  "

# Row 23 - Uniswap
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.266"
$ws.Cells.Item(23, 5).Value = "  +0.65%  "

# Row 24 - Chainlink
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "6.209"
$ws.Cells.Item(24, 5).Value = "  +0.66%  "

# Row 27 - EthereumClassic
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "18.90"
$ws.Cells.Item(27, 5).Value = "  -0.18%  "

# Row 28 - LidoDAOToken
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.919"
$ws.Cells.Item(28, 5).Value = "  -2.10%  "

# Row 29 - Stellar
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.09870"
$ws.Cells.Item(29, 5).Value = "  -1.00%  "

# Row 30 - Toncoin
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.339"
$ws.Cells.Item(30, 5).Value = "  -2.33%  "

# Row 31 - PancakeSwap
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.504"
$ws.Cells.Item(31, 5).Value = "  -0.40%  "

# Row 32 - Filecoin
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.291"
$ws.Cells.Item(32, 5).Value = "  -0.58%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.119"
$ws.Cells.Item(33, 5).Value = "  +0.84%  "

# Row 34 - Hedera
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.04829"
$ws.Cells.Item(34, 5).Value = "  +1.85%  "

# Row 35 - ARBITRUM
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.127"
$ws.Cells.Item(35, 5).Value = "  +0.41%  "

# Row 36 - ImmutableX
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.6974"
$ws.Cells.Item(36, 5).Value = "  +0.19%  "

# Row 37 - HuobiToken
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.718"
$ws.Cells.Item(37, 5).Value = "  -0.03%  "

# Row 38 - VeChain
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.01879"
$ws.Cells.Item(38, 5).Value = "  +0.90%  "

# Row 39 - MXToken
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.761"
$ws.Cells.Item(39, 5).Value = "  +0.47%  "

# Row 40 - FraxShare
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "6.286"
$ws.Cells.Item(40, 5).Value = "  +0.22%  "

# Row 41 - Aave
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "73.57"
$ws.Cells.Item(41, 5).Value = "  +5.91%  "

# Row 42 - RenderToken
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.993"
$ws.Cells.Item(42, 5).Value = "  +4.43%  "

# Row 43 - TheSandbox
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.4247"
$ws.Cells.Item(43, 5).Value = "  +1.84%  "

# Row 44 - PaxDollar
$ws.Cells.Item(44, 5).Value = "  +0.05%  "

# Row 45 - TrustWalletToken
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.8386"
$ws.Cells.Item(45, 5).Value = "  -0.34%  "

# Row 46 - Quant
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "102.24"
$ws.Cells.Item(46, 5).Value = "  +0.46%  "

# Row 47 - EnergySwap
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "9.371"
$ws.Cells.Item(47, 5).Value = "  +1.82%  "

# Row 48 - Aptos
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "7.078"
$ws.Cells.Item(48, 5).Value = "  -0.20%  "

# Row 49 - Elrond
$ws.Cells.Item(49, 5).Value = "  +0.59%  "

# Row 50 - Maker
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "913.77"
$ws.Cells.Item(50, 5).Value = "  +0.13%  "

# Row 51 - Cronos
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.05698"
$ws.Cells.Item(51, 5).Value = "  +1.91%  "
